$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value to a cell while preserving its General/text
# appearance even when the string looks numeric (e.g. "25.00", "0.999"),
# mirroring the source data which stores these as literal text.
function Set-TextValue($range, $text) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

$ws.Range("D2").Value = "42.532.38"
$ws.Range("E2").Value = "  -2.12%  "

$ws.Range("D3").Value = "2.292.80"
$ws.Range("E3").Value = "  -0.92%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.01%  "

Set-TextValue $ws.Range("D5") "302.51"
$ws.Range("E5").Value = "  -2.46%  "

Set-TextValue $ws.Range("D6") "98.36"
$ws.Range("E6").Value = "  -6.14%  "

$ws.Range("E7").Value = "  -5.34%  "

$ws.Range("E8").Value = "  +0.05%  "

Set-TextValue $ws.Range("D9") "0.499"
$ws.Range("E9").Value = "  -5.73%  "

Set-TextValue $ws.Range("D10") "34.39"

Set-TextValue $ws.Range("D11") "0.0788"
$ws.Range("E11").Value = "  -3.06%  "

Set-TextValue $ws.Range("D12") "0.113"
$ws.Range("E12").Value = "  +0.05%  "

Set-TextValue $ws.Range("D13") "6.71"
$ws.Range("E13").Value = "  -4.23%  "

$ws.Range("D14").Value = "2.644.83"
$ws.Range("E14").Value = "  -1.00%  "

Set-TextValue $ws.Range("D15") "15.65"
$ws.Range("E15").Value = "  +3.21%  "

$ws.Range("D16").Value = "2.305.06"
$ws.Range("E16").Value = "  -0.48%  "

$ws.Range("E17").Value = "  -1.25%  "

$ws.Range("D18").Value = "42.453.65"
$ws.Range("E18").Value = "  -2.09%  "

$ws.Range("D19").Value = "0.0₃0899"
$ws.Range("E19").Value = "  -2.99%  "

Set-TextValue $ws.Range("D20") "11.50"
$ws.Range("E20").Value = "  -5.61%  "

Set-TextValue $ws.Range("D21") "6.03"

Set-TextValue $ws.Range("D22") "67.72"
$ws.Range("E22").Value = "  -0.85%  "

Set-TextValue $ws.Range("D23") "234.70"
$ws.Range("E23").Value = "  -3.27%  "

Set-TextValue $ws.Range("D24") "1.97"
$ws.Range("E24").Value = "  -3.33%  "

$ws.Range("E25").Value = "  -3.31%  "

$ws.Range("E26").Value = "  -0.01%  "

Set-TextValue $ws.Range("D27") "25.00"
$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("E28").Value = "  -4.31%  "

Set-TextValue $ws.Range("D29") "34.68"
$ws.Range("E29").Value = "  -6.48%  "

Set-TextValue $ws.Range("D30") "9.16"
$ws.Range("E30").Value = "  -5.07%  "

Set-TextValue $ws.Range("D31") "163.32"
$ws.Range("E31").Value = "  -2.03%  "

Set-TextValue $ws.Range("D32") "0.999"
$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("E33").Value = "  -5.47%  "

Set-TextValue $ws.Range("D34") "4.59"
$ws.Range("E34").Value = "  +1.49%  "

$ws.Range("E35").Value = "  -4.97%  "

$ws.Range("E36").Value = "  -4.48%  "

Set-TextValue $ws.Range("D37") "16.94"
$ws.Range("E37").Value = "  -7.52%  "

$ws.Range("E38").Value = "  -6.39%  "

$ws.Range("E39").Value = "  -4.75%  "

$ws.Range("E41").Value = "  -3.73%  "

Set-TextValue $ws.Range("D42") "2.36"
$ws.Range("E42").Value = "  -13.44%  "

$ws.Range("D43").Value = "1.974.86"
$ws.Range("E43").Value = "  -0.91%  "

Set-TextValue $ws.Range("D44") "0.0280"
$ws.Range("E44").Value = "  -4.65%  "

Set-TextValue $ws.Range("D45") "18.62"
$ws.Range("E45").Value = "  -2.16%  "

Set-TextValue $ws.Range("D46") "10.14"
$ws.Range("E46").Value = "  +1.28%  "

Set-TextValue $ws.Range("D47") "2.90"
$ws.Range("E47").Value = "  -5.72%  "

Set-TextValue $ws.Range("D48") "55.50"
$ws.Range("E48").Value = "  -2.69%  "

$ws.Range("D50").Value = "2.518.00"
$ws.Range("E50").Value = "  -0.81%  "

$ws.Range("E51").Value = "  -1.30%  "
